$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cheetah processing tracking numbers - update column C (and D where mirrored)
# with the new batch of package/shipment tracking numbers, stored as text.
# These cells also lose their previous cell style (become default/unstyled),
# matching the source edit.

$updates = @{
    "C2"  = "320018179991"
    "C3"  = "320018180002"
    "C4"  = "320018180035"
    "C5"  = "320018180057"
    "D5"  = "320018180057"
    "C6"  = "320018180090"
    "D6"  = "320018180090"
    "C7"  = "320018180127"
    "D7"  = "320018180127"
    "C8"  = "320018180150"
    "C9"  = "320018180171"
    "C10" = "320018180208"
    "C11" = "320018180220"
    "C12" = "320018180263"
    "C13" = "320018180285"
    "D13" = "320018180285"
    "C14" = "320018180311"
    "D14" = "320018180311"
    "C15" = "320018180333"
    "D15" = "320018180333"
    "C16" = "320018180366"
    "D16" = "320018180366"
    "C17" = "320018180388"
    "D17" = "320018180388"
    "C18" = "320018180425"
    "C19" = "320018180447"
    "C20" = "320018180480"
    "C21" = "320018180506"
    "C22" = "320018180539"
}

# First pass: strip the old cell formatting on every target cell so they all
# share the same (default) base style before we touch NumberFormat - this
# keeps the engine from allocating a separate combined style per old style.
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# Second pass: write the new tracking numbers as text and re-flatten the
# style (NumberFormat="@" is needed so the big numeric-looking strings are
# kept as text instead of being coerced to numbers).
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
